# Apply updated Price (D) and Volume(1h) (E) values for cryptos sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '63.666.27'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +2.74%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.477.88'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +2.49%  '
$ws.Range("E4").Value = '  +0.18%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '576.94'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.62%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '149.33'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +4.40%  '
$ws.Range("E7").Value = '  -0.15%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.542'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +2.00%  '
$ws.Range("E9").Value = '  +5.36%  '
$ws.Range("E10").Value = '  +0.73%  '
$ws.Range("E11").Value = '  +3.08%  '
$ws.Range("E12").Value = '  +4.06%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '27.49'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +5.85%  '
$ws.Range("E14").Value = '  +7.39%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.949.37'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +3.33%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '63.370.53'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.37%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.488.58'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +2.95%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '11.61'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +2.16%  '
$ws.Range("E19").Value = '  +7.22%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.26'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +3.21%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '329.50'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.68%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.998'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.18%  '
$ws.Range("E23").Value = '  +11.00%  '
$ws.Range("E24").Value = '  +1.35%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '633.87'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +14.54%  '
$ws.Range("E26").Value = '  +14.22%  '
$ws.Range("E27").Value = '  +1.39%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.599.95'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +2.57%  '
$ws.Range("E29").Value = '  +10.08%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '8.50'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +3.97%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.998'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.27%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.145'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.61%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.92'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +2.58%  '
$ws.Range("E34").Value = '  +10.51%  '
$ws.Range("E35").Value = '  +4.06%  '
$ws.Range("E36").Value = '  -0.21%  '
$ws.Range("E37").Value = '  +2.38%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '5.57'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +2.94%  '
$ws.Range("E39").Value = '  +2.81%  '
$ws.Range("E40").Value = '  +3.51%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '147.38'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -4.04%  '
$ws.Range("E42").Value = '  +19.95%  '
$ws.Range("E43").Value = '  +0.70%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '151.34'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +3.16%  '
$ws.Range("E45").Value = '  +4.09%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '21.29'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +7.70%  '
$ws.Range("E47").Value = '  +4.88%  '
$ws.Range("E48").Value = '  +3.57%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0241'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +6.22%  '
$ws.Range("E50").Value = '  +1.14%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.752'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +5.82%  '
